$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") date serial value changed from 45203 (2023-10-04)
# to 45205 (2023-10-06) for all data rows (2 through 98).
$ws.Range("C2:C98").Value = 45205
